# Add a new "add_primer_features" column as the first column (A) of the
# PCRSource sheet, shifting the existing header columns one place to the
# right (circular, assembly, input, output, type, output_name, id).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PCRSource")

# Insert a new column before column A, pushing B:H to the right.
$ws.Range("A1").EntireColumn.Insert()

$ws.Range("A1").Value = "add_primer_features"
